{"js": "// Map of old multiplication-problem text -> new text (as introduced by the\n// commit). Each value is unique in the document, so an exact text match is\n// sufficient to locate the right cell.\nconst replacements = {\n  \"474\u00d75=\": \"407\u00d74=\",\n  \"759\u00d72=\": \"140\u00d73=\",\n  \"311\u00d78=\": \"746\u00d79=\",\n  \"129\u00d77=\": \"882\u00d74=\",\n  \"608\u00d76=\": \"769\u00d72=\",\n  \"972\u00d79=\": \"243\u00d77=\",\n  \"699\u00d75=\": \"520\u00d77=\",\n  \"241\u00d73=\": \"365\u00d78=\",\n  \"765\u00d72=\": \"653\u00d72=\",\n  \"920\u00d79=\": \"987\u00d78=\",\n  \"442\u00d73=\": \"159\u00d78=\",\n  \"879\u00d75=\": \"761\u00d73=\",\n  \"733\u00d72=\": \"690\u00d75=\",\n  \"125\u00d75=\": \"647\u00d72=\",\n  \"280\u00d75=\": \"139\u00d74=\",\n  \"311\u00d79=\": \"835\u00d74=\",\n  \"119\u00d74=\": \"545\u00d76=\",\n  \"397\u00d79=\": \"765\u00d79=\",\n  \"780\u00d79=\": \"969\u00d78=\",\n  \"667\u00d73=\": \"837\u00d79=\",\n  \"201\u00d76=\": \"911\u00d74=\",\n  \"802\u00d79=\": \"704\u00d72=\",\n  \"558\u00d76=\": \"385\u00d76=\",\n  \"728\u00d75=\": \"114\u00d78=\",\n  \"225\u00d76=\": \"125\u00d75=\",\n};\n\n// Walk every table cell's paragraphs and swap the text for any paragraph\n// whose full text matches one of the \"before\" keys above.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (const table of tables.items) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    const cells = row.cells;\n    cells.load(\"items\");\n    await context.sync();\n\n    for (const cell of cells.items) {\n      const paragraphs = cell.body.paragraphs;\n      paragraphs.load(\"items/text\");\n      await context.sync();\n\n      for (const paragraph of paragraphs.items) {\n        const current = paragraph.text;\n        if (Object.prototype.hasOwnProperty.call(replacements, current)) {\n          paragraph.insertText(replacements[current], \"Replace\");\n        }\n      }\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Map of old multiplication-problem text -> new text (as introduced by the\n# commit). Each value is unique in the document, so an exact text match is\n# sufficient to locate the right cell.\n$replacements = @{\n    \"474\u00d75=\" = \"407\u00d74=\";\n    \"759\u00d72=\" = \"140\u00d73=\";\n    \"311\u00d78=\" = \"746\u00d79=\";\n    \"129\u00d77=\" = \"882\u00d74=\";\n    \"608\u00d76=\" = \"769\u00d72=\";\n    \"972\u00d79=\" = \"243\u00d77=\";\n    \"699\u00d75=\" = \"520\u00d77=\";\n    \"241\u00d73=\" = \"365\u00d78=\";\n    \"765\u00d72=\" = \"653\u00d72=\";\n    \"920\u00d79=\" = \"987\u00d78=\";\n    \"442\u00d73=\" = \"159\u00d78=\";\n    \"879\u00d75=\" = \"761\u00d73=\";\n    \"733\u00d72=\" = \"690\u00d75=\";\n    \"125\u00d75=\" = \"647\u00d72=\";\n    \"280\u00d75=\" = \"139\u00d74=\";\n    \"311\u00d79=\" = \"835\u00d74=\";\n    \"119\u00d74=\" = \"545\u00d76=\";\n    \"397\u00d79=\" = \"765\u00d79=\";\n    \"780\u00d79=\" = \"969\u00d78=\";\n    \"667\u00d73=\" = \"837\u00d79=\";\n    \"201\u00d76=\" = \"911\u00d74=\";\n    \"802\u00d79=\" = \"704\u00d72=\";\n    \"558\u00d76=\" = \"385\u00d76=\";\n    \"728\u00d75=\" = \"114\u00d78=\";\n    \"225\u00d76=\" = \"125\u00d75=\";\n}\n\n# Walk every cell of every table and swap the text for any cell whose\n# (trimmed) text matches one of the \"before\" keys above.\nforeach ($t in $d.Tables) {\n    for ($r = 1; $r -le $t.Rows.Count; $r++) {\n        for ($c = 1; $c -le $t.Columns.Count; $c++) {\n            $cell = $t.Cell($r, $c)\n            $raw = $cell.Range.Text\n            # Cell.Range.Text includes the trailing cell-mark (CR + BEL); strip it.\n            $current = $raw.TrimEnd([char]13, [char]7)\n            if ($replacements.ContainsKey($current)) {\n                $cell.Range.Text = $replacements[$current]\n            }\n        }\n    }\n}\n"}
